# Add a new customer row ("寰延有限公司" / tax id 28802261) into 工作表1
# at row 251, pushing the existing rows 251..388 down to 252..389.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert a brand-new row at 251 (copies formatting from the row above,
#     i.e. the s="7"/s="7"/s="4" styles already used by similar rows). ---
$ws.Rows.Item(251).Insert()

# Match the row height used by the other single "note" rows in this sheet.
$ws.Rows.Item(251).RowHeight = 78

# Fill in C/B/A in that order so the new shared-string entries land in the
# same sequence as the target workbook (long note text, then the tax id,
# then the company name).
$ws.Cells.Item(251, 3).Value = "28802261 寰延有限公司`r`n新北市新店區德安里安康路２段２８６號（２樓）`r`n（未向國際貿易署登記出進口廠商資料者，出口金額限制美金兩萬以下，且通關必驗，若金額超過美金兩萬需檢附輸出許可證才可出口）"
$ws.Cells.Item(251, 2).Value = "28802261"
$ws.Cells.Item(251, 1).Value = "寰延有限公司"

# --- Keep the filter / named range in sync with the new, longer table. ---
$ws.Range("A1:C354").AutoFilter()

$names = $wb.Names
for ($i = 1; $i -le $names.Count(); $i++) {
    $n = $names.Item($i)
    if ($n.Name() -like "*FilterDatabase*") {
        $n.RefersTo = "=工作表1!`$A`$1:`$C`$354"
    }
}

# --- Restore the view: active cell/selection moves to the new last row. ---
$ws.Range("A372").Select()
$ws.Range("B389").Select()
